# COVID-19 Bangladesh DataSheet - "Dataset update 10 april" (last update 10 april)
#
# Appends one new daily row (row 35, date serial 43931 = 2020-04-10) to each
# of the three sheets (Confirmed, Recoverd, Death), continuing the existing
# running-total pattern:
#   - Column A: the date (copied number format from the row above)
#   - Column B: running total (a SUM formula on Confirmed/Death, a plain
#     value on Recoverd - matches the pre-existing pattern per sheet)
#   - Column C: the day's new count

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Confirmed
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Confirmed")
$ws1.Range("A34:C34").Copy()
$ws1.Range("A35:C35").PasteSpecial(-4122)   # xlPasteFormats - inherit styles from row 34
$ws1.Range("A35").Value = 43931
$ws1.Range("B35").Formula = "=SUM(B34+C35)"
$ws1.Range("C35").Value = 94

# ---------------------------------------------------------------------
# Recoverd
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Recoverd")
$ws2.Range("A34:C34").Copy()
$ws2.Range("A35:C35").PasteSpecial(-4122)
$ws2.Range("A35").Value = 43931
$ws2.Range("B35").Value = 33
$ws2.Range("C35").Value = 0

# ---------------------------------------------------------------------
# Death
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Death")
$ws3.Range("A34:C34").Copy()
$ws3.Range("A35:C35").PasteSpecial(-4122)
$ws3.Range("A35").Value = 43931
$ws3.Range("B35").Formula = "=SUM(B34+C35)"
$ws3.Range("C35").Value = 6

# ---------------------------------------------------------------------
# View/selection state to match the saved workbook: Recoverd used to be the
# active tab with C35 selected; now Confirmed is active (K22 selected),
# Recoverd just has B35 selected and Death has D36 selected.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("B35").Select()

$ws3.Activate()
$ws3.Range("D36").Select()

$ws1.Activate()
$ws1.Range("K22").Select()
